# Bug fix in factory connections:
# The CO2 emitted from the BF (blast furnace) outflow used a single generic
# "CO2__emitted" flow. Split this into two separate flows so fossil and
# biogenic CO2 can be tracked/connected separately into CO2 capture:
#   - rename the existing connection's product to "CO2__fossil"
#   - add a duplicate connection row for "CO2__bio"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connections")

# Row 20 previously carried the product "CO2__emitted" -> rename to CO2__fossil
$ws.Range("E20").Value = "CO2__fossil"

# Insert a new row 21 that duplicates row 20 (same chain/unit/flow types and
# styles), then change its product to CO2__bio. This mirrors the fossil CO2
# connection for the biogenic CO2 stream.
$ws.Rows.Item(20).Copy()
$ws.Rows.Item(21).Insert()
$ws.Range("E21").Value = "CO2__bio"

$ws.Range("J33").Select()
